$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "284×7=1988" "665×2=1330"
Replace-Text "551×4=2204" "169×9=1521"
Replace-Text "378×4=1512" "353×3=1059"
Replace-Text "168×2=336" "634×5=3170"
Replace-Text "313×7=2191" "240×8=1920"
Replace-Text "976×3=2928" "518×5=2590"
Replace-Text "895×9=8055" "212×6=1272"
Replace-Text "877×9=7893" "669×6=4014"
Replace-Text "582×7=4074" "786×3=2358"
Replace-Text "662×8=5296" "727×2=1454"
Replace-Text "412×7=2884" "708×9=6372"
Replace-Text "716×7=5012" "178×7=1246"
Replace-Text "239×7=1673" "701×4=2804"
Replace-Text "842×8=6736" "546×3=1638"
Replace-Text "313×5=1565" "354×7=2478"
Replace-Text "589×3=1767" "139×5=695"
Replace-Text "125×9=1125" "530×4=2120"
Replace-Text "442×9=3978" "283×2=566"
Replace-Text "466×5=2330" "573×2=1146"
Replace-Text "869×4=3476" "210×2=420"
Replace-Text "648×4=2592" "868×5=4340"
Replace-Text "616×3=1848" "900×9=8100"
Replace-Text "134×6=804" "796×6=4776"
Replace-Text "983×6=5898" "715×2=1430"
Replace-Text "360×9=3240" "123×7=861"
